# Scheduled runner update: refresh market-price columns (H-N) on a handful
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4200
$ws.Range("J43").Value = 4333.6665
$ws.Range("L43").Value = 4333.6665
$ws.Range("N43").Value = -4471.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 390.4
$ws.Range("I53").Value = 152.25
$ws.Range("J53").Value = 477
$ws.Range("K53").Value = 152.25
$ws.Range("L53").Value = 477
$ws.Range("M53").Value = 484.75
$ws.Range("N53").Value = -1751

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1982.1428
$ws.Range("I82").Value = 1982.1428
$ws.Range("K82").Value = 5946.428400000001
$ws.Range("M82").Value = -5540.428400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 1982.1428
$ws.Range("I85").Value = 1982.1428
$ws.Range("K85").Value = 5946.428400000001
$ws.Range("M85").Value = -4542.428400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14625.357
$ws.Range("I132").Value = 2160.7144
$ws.Range("K132").Value = 6482.1432
$ws.Range("M132").Value = -3952.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2388.158
$ws.Range("I137").Value = 2159.8064
$ws.Range("K137").Value = 6479.4192
$ws.Range("M137").Value = -3929.4192

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3513.9575
$ws.Range("I138").Value = 2633.75
$ws.Range("J138").Value = 3694.513
$ws.Range("K138").Value = 7901.25
$ws.Range("L138").Value = 11083.539
$ws.Range("M138").Value = -2761.25
$ws.Range("N138").Value = -21363.539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 465
$ws.Range("J17").Value = 465
$ws.Range("L17").Value = 465
$ws.Range("N17").Value = -811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2177.9412
$ws.Range("J61").Value = 3831.6667
$ws.Range("L61").Value = 3831.6667
$ws.Range("N61").Value = -4255.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3906.5908
$ws.Range("I110").Value = 3709.8235
$ws.Range("K110").Value = 3709.8235
$ws.Range("M110").Value = -1664.8235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2067.7917
$ws.Range("I132").Value = 2351.6875
$ws.Range("K132").Value = 7055.0625
$ws.Range("M132").Value = -4525.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2177.9412
$ws.Range("J136").Value = 3831.6667
$ws.Range("L136").Value = 11495.0001
$ws.Range("N136").Value = -16595.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3550.7297
$ws.Range("I107").Value = 3203.742
$ws.Range("J107").Value = 5343.5
$ws.Range("K107").Value = 3203.742
$ws.Range("L107").Value = 5343.5
$ws.Range("M107").Value = -1283.742
$ws.Range("N107").Value = -9183.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3619.5833
$ws.Range("I134").Value = 3063.1428
$ws.Range("K134").Value = 9189.4284
$ws.Range("M134").Value = -6654.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1504.2593
$ws.Range("I31").Value = 1424.6
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1424.6
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1129.6
$ws.Range("N31").Value = -3090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1504.2593
$ws.Range("I34").Value = 1424.6
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1424.6
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1222.6
$ws.Range("N34").Value = -2904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3683.85
$ws.Range("J99").Value = 3362.3333
$ws.Range("L99").Value = 3362.3333
$ws.Range("N99").Value = -6358.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3683.85
$ws.Range("J126").Value = 3362.3333
$ws.Range("L126").Value = 10086.9999
$ws.Range("N126").Value = -15026.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3006.4443
$ws.Range("I132").Value = 2379.4285
$ws.Range("J132").Value = 3681.6924
$ws.Range("K132").Value = 7138.2855
$ws.Range("L132").Value = 11045.0772
$ws.Range("M132").Value = -4608.2855
$ws.Range("N132").Value = -16105.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7155.7144
$ws.Range("I56").Value = 7155.7144
$ws.Range("K56").Value = 7155.7144
$ws.Range("M56").Value = -6625.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5768.846
$ws.Range("I80").Value = 6623.75
$ws.Range("K80").Value = 19871.25
$ws.Range("M80").Value = -18935.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5768.846
$ws.Range("I83").Value = 6623.75
$ws.Range("K83").Value = 59613.75
$ws.Range("M83").Value = -54933.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3274.8845
$ws.Range("J122").Value = 8541.777
$ws.Range("L122").Value = 76875.993
$ws.Range("N122").Value = -81775.993

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 100003550
$ws.Range("I141").Value = 100003550
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 300010650
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -300005470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 8849108
$ws.Range("I24").Value = 11500840
$ws.Range("J24").Value = 10000
$ws.Range("K24").Value = 11500840
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = -11500667
$ws.Range("N24").Value = -10346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8499.5
$ws.Range("I41").Value = 1999
$ws.Range("K41").Value = 1999
$ws.Range("M41").Value = -1644

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 96701.5
$ws.Range("I70").Value = 161980.58
$ws.Range("K70").Value = 161980.58
$ws.Range("M70").Value = -161710.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 96701.5
$ws.Range("I73").Value = 161980.58
$ws.Range("K73").Value = 161980.58
$ws.Range("M73").Value = -161044.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 651.3182
$ws.Range("I107").Value = 547.58826
$ws.Range("J107").Value = 1004
$ws.Range("K107").Value = 547.58826
$ws.Range("L107").Value = 1004
$ws.Range("M107").Value = 1372.41174
$ws.Range("N107").Value = -4844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 39496.25
$ws.Range("J123").Value = 39496.25
$ws.Range("L123").Value = 39496.25
$ws.Range("N123").Value = -44396.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 83442
$ws.Range("J134").Value = 83442
$ws.Range("L134").Value = 250326
$ws.Range("N134").Value = -255396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 14706374
$ws.Range("J16").Value = 625
$ws.Range("L16").Value = 625
$ws.Range("N16").Value = -965

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6061598.5
$ws.Range("I22").Value = 10101455
$ws.Range("K22").Value = 10101455
$ws.Range("M22").Value = -10101160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 6061598.5
$ws.Range("I27").Value = 10101455
$ws.Range("K27").Value = 10101455
$ws.Range("M27").Value = -10101348

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1026.3334
$ws.Range("I55").Value = 2969
$ws.Range("J55").Value = 55
$ws.Range("K55").Value = 2969
$ws.Range("L55").Value = 55
$ws.Range("M55").Value = -2796
$ws.Range("N55").Value = -401

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 84130
$ws.Range("J86").Value = 84130
$ws.Range("L86").Value = 84130
$ws.Range("N86").Value = -86376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 84130
$ws.Range("J89").Value = 84130
$ws.Range("L89").Value = 420650
$ws.Range("N89").Value = -431882
